$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data cells store text-typed values ("1", "2", "3", "4" are stored as
# shared-string text, not numbers), so plain .Value assignments of numeric-
# looking strings would get auto-coerced to numeric cells. Instead, copy the
# new contents from existing text cells that already hold the same values,
# which preserves the text cell type without touching cell styling.

# Build new row 7 first (before row 2 is changed below), reusing:
#  - A2:D2 which already holds "a","1","1","1"
#  - E5:H5 which already holds "c","3","3","3"
$ws.Range("A2:D2").Copy()
$ws.Range("A7:D7").PasteSpecial()
$ws.Range("E5:H5").Copy()
$ws.Range("E7:H7").PasteSpecial()

# Update row 2: columns E:H change from "c","3","3","3" to "b","2","2","2",
# reusing E4:H4 which already holds "b","2","2","2".
$ws.Range("E4:H4").Copy()
$ws.Range("E2:H2").PasteSpecial()

$excel.CutCopyMode = $false
